$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.334.64'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.885.13'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.35%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.93'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.83%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4693'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.63%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06613'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.63'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +10.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07774'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '98.05'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -3.34%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.883.70'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.086'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6772'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '285.04'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +8.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.316.55'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9997'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.64'
$ws.Range("D19").ClearFormats()
$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.129.97'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.22%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.415'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.77%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007287'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.193'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.394'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.31'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.80%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.24'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.79%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.993'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.372'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09682'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.386'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -7.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.482'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.136'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.49%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04676'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7066'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.45%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.100'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.42%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.717'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01871'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.585'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +5.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.524'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.86%  '
$ws.Range("E41").Value = '  -3.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.970'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8664'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.63%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '103.09'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.41%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4188'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '982.37'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +8.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.292'
$ws.Range("D48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.188'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +4.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.94'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1146'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -4.47%  '
